# Tool Practice Learning Package -- apply the edits described by the
# commit ("Minor rewording" / "Minor fixes").
#
# Strategy: text-level changes are applied with Find/Replace (wildcards
# off) on $d.Content, which matches across run boundaries and merges the
# runs it touches -- that mirrors the run re-splits seen in the source
# diff (the actual rendered/stored text is what matters). The numbering
# (bullet glyph) swap between the two numbered lists is applied through
# ListFormat/ListTemplate/ListLevel, which is the only lever the object
# model exposes for that kind of change.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# --- Goals paragraph -------------------------------------------------
Replace-Text "make sure that our target audience" "ensure that our target audience"

# --- Target Audience paragraph ---------------------------------------
Replace-Text "Our primary target audience will be those within a development team." "Our primary target audience will be the individuals within a development team."

# --- Learning Plan intro paragraph ------------------------------------
Replace-Text "the learning package will be divided up into four sessions" "the learning package will be divided into four sessions"
Replace-Text "to keep the student from losing focus." "to keep the student from staying focused."

# --- Session 1 ---------------------------------------------------------
Replace-Text "The training will then aim to familiarize the student with the Trello environment." "The training will then aim to get the student be familiar with the Trello environment."
Replace-Text "As an account will be needed to carry on in the training, all students will be checked to see if they have successfully created their account" "As an account will be needed to carry on in the training, all each will be required to undergo an assessment to see if they have successfully created their account"

# --- Session 2 ---------------------------------------------------------
Replace-Text "Student will be placed into groups of two or more and will practice" "Student will be placed into groups of two or more people and will practice"
Replace-Text "They will then be taught about discussion tools and how collaboration works" "They will then be taught about the discussion tools and how collaboration works"
Replace-Text "If time permits, student will be taught in depth functions of cards such as labels, due date, checklist, attachments etc." "If time permits, student will be taught in-depth functions of cards such as labels, due date, checklist, attachments and so on."

# --- Session 4 ---------------------------------------------------------
Replace-Text "In this session, students will be taught about advanced functions of Trello boards and how to effectively use them to make the most of their projects." "In this session, students will be taught about the advanced functions of Trello boards and as to how they could efficiently use them to make the most of their projects."

# --- Learning Activity ---------------------------------------------------------
Replace-Text "d with at least 5 user stories with labels and tasks. Students will also be tasked to assign members" "d with at least 5 user stories comprising of labels and tasks. Students will also be obligated to assign members"

# --- Numbering: swap the bullet glyph "families" used by the two
# numbered lists in the document (numId 1 -- the feature list under
# "Here are some of the most important features..." -- and numId 2 --
# the Session 1-4 list). In the target, numId 1 ends up using the plain
# dash glyph and numId 2 ends up using the rotating bullet/circle/square
# glyphs; right now it is the other way round.
$featurePara = $null
$sessionPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($featurePara -eq $null -and $t -like "*Real-time collaboration*") {
        $featurePara = $p
    }
    if ($sessionPara -eq $null -and $t -like "*Session 1*") {
        $sessionPara = $p
    }
    if ($featurePara -ne $null -and $sessionPara -ne $null) {
        break
    }
}

$featureTemplate = $featurePara.Range.ListFormat.ListTemplate
$sessionTemplate = $sessionPara.Range.ListFormat.ListTemplate

$featureFormats = @()
$sessionFormats = @()
for ($lvl = 1; $lvl -le 9; $lvl++) {
    $featureFormats += $featureTemplate.ListLevels.Item($lvl).NumberFormat
    $sessionFormats += $sessionTemplate.ListLevels.Item($lvl).NumberFormat
}
for ($lvl = 1; $lvl -le 9; $lvl++) {
    $featureTemplate.ListLevels.Item($lvl).NumberFormat = $sessionFormats[$lvl - 1]
    $sessionTemplate.ListLevels.Item($lvl).NumberFormat = $featureFormats[$lvl - 1]
}
